$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("RESEARCH GROUP / DESIGN TEAM ", $true, $false, $false, $false, $false, $true, 1, $false, "RESEARCH GROUP ", 2)
